$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 4826
$ws.Range("L3").Value = 5192
$ws.Range("K4").Value = 1297
$ws.Range("L4").Value = 1273
$ws.Range("L6").Value = 4391
$ws.Range("K7").Value = 20033
$ws.Range("L7").Value = 15986

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L6").Value = 121
$ws.Range("L7").Value = 524
$ws.Range("L8").Value = 1064
$ws.Range("L15").Value = 117
$ws.Range("L19").Value = 438
$ws.Range("L20").Value = 401
$ws.Range("L22").Value = 47
$ws.Range("L27").Value = 143
$ws.Range("L29").Value = 875
$ws.Range("L31").Value = 163
$ws.Range("L33").Value = 732
$ws.Range("L34").Value = 93
$ws.Range("L36").Value = 208
$ws.Range("L37").Value = 600
$ws.Range("L42").Value = 519
$ws.Range("L48").Value = 207
$ws.Range("L50").Value = 79
$ws.Range("L51").Value = 201
$ws.Range("L52").Value = 323
$ws.Range("L53").Value = 181
$ws.Range("L54").Value = 335
$ws.Range("J63").Value = 124
$ws.Range("K63").Value = 150
$ws.Range("L63").Value = 44
$ws.Range("L65").Value = 312
$ws.Range("L67").Value = 551
$ws.Range("L76").Value = 247
$ws.Range("L83").Value = 354
$ws.Range("L84").Value = 155
$ws.Range("L85").Value = 819
$ws.Range("L90").Value = 161
$ws.Range("J95").Value = 307
$ws.Range("L100").Value = 26
$ws.Range("K101").Value = 20033
$ws.Range("L101").Value = 15986

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 178
$ws.Range("L3").Value = 175
$ws.Range("L7").Value = 524

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 246
$ws.Range("L6").Value = 170
$ws.Range("L7").Value = 819

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L2").Value = 105
$ws.Range("L6").Value = 88
$ws.Range("L7").Value = 323

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L2").Value = 54
$ws.Range("L7").Value = 181

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 310
$ws.Range("L3").Value = 358
$ws.Range("L6").Value = 277
$ws.Range("L7").Value = 1064

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L3").Value = 142
$ws.Range("L7").Value = 354

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 200
$ws.Range("L3").Value = 251
$ws.Range("L6").Value = 222
$ws.Range("L7").Value = 732

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J4").Value = 15
$ws.Range("J7").Value = 307

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 179
$ws.Range("L3").Value = 204
$ws.Range("L4").Value = 34
$ws.Range("L6").Value = 165
$ws.Range("L7").Value = 600

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L4").Value = 15
$ws.Range("L7").Value = 312

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L3").Value = 116
$ws.Range("L6").Value = 57

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("L2").Value = 64
$ws.Range("L3").Value = 42
$ws.Range("L6").Value = 45
$ws.Range("L7").Value = 163

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 161
$ws.Range("L6").Value = 126
$ws.Range("L7").Value = 551

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L6").Value = 45
$ws.Range("L7").Value = 155

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L6").Value = 163
$ws.Range("L7").Value = 335

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 262
$ws.Range("L3").Value = 330
$ws.Range("L7").Value = 875

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L3").Value = 53
$ws.Range("L7").Value = 207

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L6").Value = 124
$ws.Range("L7").Value = 438

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L2").Value = 49
$ws.Range("L6").Value = 114
$ws.Range("L7").Value = 247

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("L2").Value = 52
$ws.Range("L7").Value = 121

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L2").Value = 150
$ws.Range("L7").Value = 519

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L2").Value = 124
$ws.Range("L6").Value = 108
$ws.Range("L7").Value = 401

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L3").Value = 62
$ws.Range("L7").Value = 208

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("L6").Value = 14
$ws.Range("L7").Value = 26

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("L3").Value = 25
$ws.Range("L7").Value = 93

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L3").Value = 38
$ws.Range("L7").Value = 117

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("L2").Value = 28
$ws.Range("L7").Value = 79

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L3").Value = 42
$ws.Range("L7").Value = 143

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L2").Value = 54
$ws.Range("L7").Value = 161

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L2").Value = 58
$ws.Range("L7").Value = 201

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("L3").Value = 18
$ws.Range("L7").Value = 47
